# ppt: moved ASU logo
#
# Slide 1 ("Digital Humanities" title slide): the ASU master-logo picture is
# moved up/shrunk slightly, and the subtitle placeholder below it is
# repositioned/resized to match.
#
# NOTE: PowerPoint's COM object model reports/accepts Shape.Left/Top/Width/
# Height in points (1 pt = 12700 EMU) and stores them internally as 32-bit
# floats, while the underlying OOXML stores exact integer EMUs. A naive
# EMU/12700.0 conversion can therefore land 1 EMU off after the point value
# is rounded to float32 precision and converted back. EmuToPt searches the
# neighborhood of the naive conversion for a point value that round-trips
# to the exact target EMU.

function EmuToPt {
    param([long]$EmuTarget)

    $base = [float]($EmuTarget / 12700.0)

    $test = [long]([math]::Floor([double]$base * 12700.0))
    if ($test -eq $EmuTarget) {
        return [double]$base
    }

    $step = 0.0000002
    for ($i = 1; $i -le 200000; $i++) {
        foreach ($sign in 1, -1) {
            $cand = [float]($base + ($sign * $i * $step))
            $back = [double]$cand * 12700.0
            $emu = [long]([math]::Floor($back))
            if ($emu -eq $EmuTarget) {
                return [double]$cand
            }
        }
    }

    # Fall back to the naive conversion if no exact float32 match is found.
    return [double]$base
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Picture 3 - ASU master logo
$logo = $s.Shapes.Item(3)
$logo.Left   = EmuToPt 2737834
$logo.Top    = EmuToPt 1880088
$logo.Width  = EmuToPt 8072867
$logo.Height = EmuToPt 1777512

# Subtitle 2 - repositioned/resized under the relocated logo
$subtitle = $s.Shapes.Item(2)
$subtitle.Left   = EmuToPt 2737834
$subtitle.Top    = EmuToPt 4794963
$subtitle.Width  = EmuToPt 8067064
$subtitle.Height = EmuToPt 1126283
